$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 23 de Octubre de 2020 a las 10:48"

# Update country names (column A) where the sort order changed
$ws.Range("A31").Value = "Polonia"
$ws.Range("A32").Value = "Chequia"

$ws.Range("A75").Value = "Tunez"
$ws.Range("A76").Value = "Kenia"
$ws.Range("A77").Value = "Jordania"

$ws.Range("A113").Value = "Lituania"
$ws.Range("A114").Value = "Haiti"
$ws.Range("A115").Value = "Gabon"

$ws.Range("A142").Value = "Estonia"
$ws.Range("A143").Value = "Islandia"

$ws.Range("A216").Value = "Islas Malvinas"
$ws.Range("A217").Value = "Montserrat"

# Update numeric data (columns B-H) for all affected rows
function Set-Row($r, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = $c
    $ws.Cells.Item($r, 4).Value = $d
    $ws.Cells.Item($r, 5).Value = $e
    $ws.Cells.Item($r, 6).Value = $f
    $ws.Cells.Item($r, 7).Value = $g
    $ws.Cells.Item($r, 8).Value = $h
}

Set-Row 7   1480646 17340 1119251 335870 0 283 25525
Set-Row 23  365799  1923  312691  46193  0 132 6915
Set-Row 31  228318  13632 105092  119054 0 153 4172
Set-Row 32  223065  0     87225   133995 0 0   1845
Set-Row 65  57951   10    57829   94     0 0   28
Set-Row 75  47214   1322  5032    41398  0 44  784
Set-Row 76  47212   0     33050   13292  0 0   870
Set-Row 77  46441   0     7340    38620  0 0   481
Set-Row 86  31717   1867  22910   8394   0 7   413
Set-Row 113 9104    442   3978    5000   0 1   126
Set-Row 114 9007    0     7311    1465   0 0   231
Set-Row 115 8901    0     8479    368    0 0   54
Set-Row 142 4300    53    3418    811    0 0   71
Set-Row 143 4268    0     3098    1159   0 0   11
Set-Row 216 13      0     13      0      0 0   0
Set-Row 217 13      0     12      0      0 0   1
Set-Row 221 1       0     1       0      0 0   0
